# Adds a Wingdings checkmark symbol (the same one already used on rows like
# "Winning Condition") to the end of the "Create Main Menu" and
# "Manage Scene Transitions" progress-table rows.
#
# Range.InsertSymbol is a no-op in this host, so the checkmark run is
# injected as raw OOXML via Range.InsertXML instead. InsertXML replaces the
# addressed range's content, so the range is expanded to cover the whole
# paragraph (minus its trailing paragraph mark) and the XML payload re-states
# that paragraph's existing runs followed by the new <w:sym> run; wrapping
# the runs in a bare <w:p> (no attributes/pPr) makes the host keep the
# paragraph's own identity/pPr and only swap in the supplied run content.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Add-ChecklistCheckmark($paragraphStartText, $existingRunsXml) {
    $rng = $d.Content
    $rng.Find.Execute($paragraphStartText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Expand(4) | Out-Null          # wdParagraph
    $rng.MoveEnd(1, -1) | Out-Null     # drop the trailing paragraph mark

    $checkRun = "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:sym w:font=`"Wingdings`" w:char=`"F0FC`"/></w:r>"
    $xml = "<w:p $wNs>$existingRunsXml$checkRun</w:p>"
    $rng.InsertXML($xml)
}

$createMainMenuRuns =
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>Create Main Menu</w:t></w:r>" +
    "<w:r w:rsidR=`"007A78C8`"><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:tab/></w:r>" +
    "<w:r w:rsidR=`"007A78C8`"><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:tab/><w:t>-</w:t></w:r>" +
    "<w:r w:rsidR=`"007A78C8`"><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:tab/></w:r>"

$manageSceneTransitionsRuns =
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>Manage Scene Transitions</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:tab/></w:r>" +
    "<w:r w:rsidR=`"007A78C8`"><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>-</w:t></w:r>" +
    "<w:r w:rsidR=`"007A78C8`"><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr><w:tab/></w:r>"

Add-ChecklistCheckmark "Create Main Menu" $createMainMenuRuns
Add-ChecklistCheckmark "Manage Scene Transitions" $manageSceneTransitionsRuns

Write-Output "done"
